# Auto-generated Excel COM-interop script
# Applies updated profit-calculation values to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# (per scheduled-runner refresh of market-price-derived columns H-N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 710.9091
$ws.Range("I8").Value = 202.22223
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 606.66669
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = -467.66669
$ws.Range("N8").Value = -9278

# Row 9
$ws.Range("H9").Value = 135.2
$ws.Range("I9").Value = 160.16667
$ws.Range("J9").Value = 97.75
$ws.Range("K9").Value = 160.16667
$ws.Range("L9").Value = 97.75
$ws.Range("M9").Value = 8.833329999999989
$ws.Range("N9").Value = -435.75

# Row 53
$ws.Range("H53").Value = 190.71428
$ws.Range("I53").Value = 206.66667
$ws.Range("K53").Value = 206.66667
$ws.Range("M53").Value = 430.33333

# Row 70
$ws.Range("H70").Value = 4949.6665
$ws.Range("I70").Value = 1600
$ws.Range("K70").Value = 4800
$ws.Range("M70").Value = -4530

# Row 73
$ws.Range("H73").Value = 4949.6665
$ws.Range("I73").Value = 1600
$ws.Range("K73").Value = 4800
$ws.Range("M73").Value = -3864

# Row 76
$ws.Range("H76").Value = 4499.6665
$ws.Range("I76").Value = 4499.5
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 4499.5
$ws.Range("L76").Value = 4500
$ws.Range("M76").Value = -4184.5
$ws.Range("N76").Value = -5130

# Row 79
$ws.Range("H79").Value = 4499.6665
$ws.Range("I79").Value = 4499.5
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 4499.5
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = -3407.5
$ws.Range("N79").Value = -6684

# Row 88
$ws.Range("H88").Value = 1103.625
$ws.Range("I88").Value = 1250
$ws.Range("J88").Value = 1054.8334
$ws.Range("K88").Value = 1250
$ws.Range("L88").Value = 1054.8334
$ws.Range("M88").Value = -844
$ws.Range("N88").Value = -1866.8334

# Row 91
$ws.Range("H91").Value = 1103.625
$ws.Range("I91").Value = 1250
$ws.Range("J91").Value = 1054.8334
$ws.Range("K91").Value = 1250
$ws.Range("L91").Value = 1054.8334
$ws.Range("M91").Value = 154
$ws.Range("N91").Value = -3862.8334

# Row 116
$ws.Range("H116").Value = 2137.8235
$ws.Range("I116").Value = 2043.3077
$ws.Range("J116").Value = 2445
$ws.Range("K116").Value = 2043.3077
$ws.Range("L116").Value = 2445
$ws.Range("M116").Value = 1398.6923
$ws.Range("N116").Value = -9329

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 908.46155
$ws.Range("I2").Value = 908.46155
$ws.Range("K2").Value = 908.46155
$ws.Range("M2").Value = -795.46155

# Row 32
$ws.Range("H32").Value = 9681.223
$ws.Range("I32").Value = 10266.375
$ws.Range("K32").Value = 10266.375
$ws.Range("M32").Value = -9979.375

# Row 39
$ws.Range("H39").Value = 4284.1665
$ws.Range("I39").Value = 3141
$ws.Range("K39").Value = 3141
$ws.Range("M39").Value = -2621

# Row 61
$ws.Range("H61").Value = 6249.1665
$ws.Range("I61").Value = 4331.6665
$ws.Range("J61").Value = 8166.6665
$ws.Range("K61").Value = 4331.6665
$ws.Range("L61").Value = 8166.6665
$ws.Range("M61").Value = -4119.6665
$ws.Range("N61").Value = -8590.666499999999

# Row 116
$ws.Range("H116").Value = 908.46155
$ws.Range("I116").Value = 908.46155
$ws.Range("K116").Value = 908.46155
$ws.Range("M116").Value = 1385.53845

# Row 136
$ws.Range("H136").Value = 6249.1665
$ws.Range("I136").Value = 4331.6665
$ws.Range("J136").Value = 8166.6665
$ws.Range("K136").Value = 12994.9995
$ws.Range("L136").Value = 24499.9995
$ws.Range("M136").Value = -10444.9995
$ws.Range("N136").Value = -29599.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 908.46155
$ws.Range("I3").Value = 908.46155
$ws.Range("K3").Value = 908.46155
$ws.Range("M3").Value = -794.46155

# Row 54
$ws.Range("H54").Value = 1810.8
$ws.Range("I54").Value = 763.5
$ws.Range("K54").Value = 763.5
$ws.Range("M54").Value = -279.5

# Row 86
$ws.Range("H86").Value = 6329.5264
$ws.Range("I86").Value = 5163.6665
$ws.Range("K86").Value = 5163.6665
$ws.Range("M86").Value = -4040.6665

# Row 89
$ws.Range("H89").Value = 6329.5264
$ws.Range("I89").Value = 5163.6665
$ws.Range("K89").Value = 25818.3325
$ws.Range("M89").Value = -20202.3325

# Row 107
$ws.Range("H107").Value = 3614.9
$ws.Range("I107").Value = 2392.8125
$ws.Range("K107").Value = 2392.8125
$ws.Range("M107").Value = -472.8125

$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 304.5
$ws.Range("I13").Value = 304.5
$ws.Range("K13").Value = 304.5
$ws.Range("M13").Value = -165.5

# Row 22
$ws.Range("H22").Value = 1900.091
$ws.Range("I22").Value = 1737.375
$ws.Range("K22").Value = 1737.375
$ws.Range("M22").Value = -1387.375

# Row 99
$ws.Range("H99").Value = 2379.7
$ws.Range("I99").Value = 1994
$ws.Range("K99").Value = 1994
$ws.Range("M99").Value = -496

# Row 125
$ws.Range("H125").Value = 90000
$ws.Range("J125").Value = 90000
$ws.Range("L125").Value = 90000
$ws.Range("N125").Value = -94920

# Row 126
$ws.Range("H126").Value = 2379.7
$ws.Range("I126").Value = 1994
$ws.Range("K126").Value = 5982
$ws.Range("M126").Value = -3512

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 617
$ws.Range("I34").Value = 107.333336
$ws.Range("J34").Value = 871.8333
$ws.Range("K34").Value = 322.000008
$ws.Range("L34").Value = 2615.4999
$ws.Range("M34").Value = -238.000008
$ws.Range("N34").Value = -2783.4999

# Row 47
$ws.Range("H47").Value = 91.666664
$ws.Range("I47").Value = 91.666664
$ws.Range("K47").Value = 274.999992
$ws.Range("M47").Value = 156.000008

# Row 57
$ws.Range("H57").Value = 1363.1578
$ws.Range("J57").Value = 1363.1578
$ws.Range("L57").Value = 4089.4734
$ws.Range("N57").Value = -5207.4734

# Row 86
$ws.Range("H86").Value = 390.55554
$ws.Range("J86").Value = 420
$ws.Range("L86").Value = 1260
$ws.Range("N86").Value = -3632

# Row 89
$ws.Range("H89").Value = 390.55554
$ws.Range("J89").Value = 420
$ws.Range("L89").Value = 3780
$ws.Range("N89").Value = -15636

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 3260
$ws.Range("I70").Value = 3307.5557
$ws.Range("J70").Value = 3153
$ws.Range("K70").Value = 3307.5557
$ws.Range("L70").Value = 3153
$ws.Range("M70").Value = -3037.5557
$ws.Range("N70").Value = -3693

# Row 73
$ws.Range("H73").Value = 3260
$ws.Range("I73").Value = 3307.5557
$ws.Range("J73").Value = 3153
$ws.Range("K73").Value = 3307.5557
$ws.Range("L73").Value = 3153
$ws.Range("M73").Value = -2371.5557
$ws.Range("N73").Value = -5025

# Row 122
$ws.Range("H122").Value = 2335.75
$ws.Range("I122").Value = 1558.4286
$ws.Range("J122").Value = 7777
$ws.Range("K122").Value = 4675.2858
$ws.Range("L122").Value = 23331
$ws.Range("M122").Value = -2225.2858
$ws.Range("N122").Value = -28231

$ws = $wb.Worksheets.Item("LTW")
# Row 104
$ws.Range("H104").Value = 18370
$ws.Range("J104").Value = 18370
$ws.Range("L104").Value = 18370
$ws.Range("N104").Value = -25358

$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 21334.666
$ws.Range("I3").Value = 3500
$ws.Range("J3").Value = 57004
$ws.Range("K3").Value = 3500
$ws.Range("L3").Value = 57004
$ws.Range("M3").Value = -3386
$ws.Range("N3").Value = -57232
